$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.022286283269705
$ws.Range("D2").Value = 1.027019936693963
$ws.Range("E2").Value = 1.047095805259104
$ws.Range("F2").Value = 1.051020231627633
$ws.Range("I2").Value = 1.028335497973462
$ws.Range("J2").Value = 1.027472974458237
$ws.Range("K2").Value = 1.029840884956849
$ws.Range("L2").Value = 1.04985940014644
$ws.Range("M2").Value = 1.053772892728296
$ws.Range("N2").Value = 1.028932103934765

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.023260816035162
$ws.Range("D3").Value = 1.027724678571712
$ws.Range("E3").Value = 1.048387890737887
$ws.Range("F3").Value = 1.052390017255731
$ws.Range("I3").Value = 1.028460813560444
$ws.Range("J3").Value = 1.028085360982689
$ws.Range("K3").Value = 1.030353542333091
$ws.Range("L3").Value = 1.05096201487684
$ws.Range("M3").Value = 1.054953807404114
$ws.Range("N3").Value = 1.029545360118324

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.023891330244021
$ws.Range("D4").Value = 1.028180371058632
$ws.Range("E4").Value = 1.049224590461323
$ws.Range("F4").Value = 1.053277048115614
$ws.Range("I4").Value = 1.028540392876332
$ws.Range("J4").Value = 1.028480968289222
$ws.Range("K4").Value = 1.030684272980284
$ws.Range("L4").Value = 1.051675565217991
$ws.Range("M4").Value = 1.055718092020004
$ws.Range("N4").Value = 1.029941529232608

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.0241563807752
$ws.Range("D5").Value = 1.028371866119829
$ws.Range("E5").Value = 1.049576492789192
$ws.Range("F5").Value = 1.053650122639757
$ws.Range("I5").Value = 1.028573486607616
$ws.Range("J5").Value = 1.028647126295204
$ws.Range("K5").Value = 1.03082307392949
$ws.Range("L5").Value = 1.051975563811129
$ws.Range("M5").Value = 1.056039435805985
$ws.Range("N5").Value = 1.030107923202018

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.024200882890864
$ws.Range("D6").Value = 1.028404014381753
$ws.Range("E6").Value = 1.049635587812918
$ws.Range("F6").Value = 1.053712773332035
$ws.Range("I6").Value = 1.028579021987506
$ws.Range("J6").Value = 1.028675015841345
$ws.Range("K6").Value = 1.030846365236
$ws.Range("L6").Value = 1.052025936236587
$ws.Range("M6").Value = 1.056093393150775
$ws.Range("N6").Value = 1.030135852354514

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.023894871931527
$ws.Range("D7").Value = 1.028182930132784
$ws.Range("E7").Value = 1.04922929199078
$ws.Range("F7").Value = 1.053282032497876
$ws.Range("I7").Value = 1.02854083649692
$ws.Range("J7").Value = 1.028483189111185
$ws.Range("K7").Value = 1.030686128582023
$ws.Range("L7").Value = 1.051679573725171
$ws.Range("M7").Value = 1.055722385677859
$ws.Range("N7").Value = 1.029943753208394

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.022615646426624
$ws.Range("D8").Value = 1.027258173822297
$ws.Range("E8").Value = 1.047532341507748
$ws.Range("F8").Value = 1.051483015879425
$ws.Range("I8").Value = 1.028378161054483
$ws.Range("J8").Value = 1.02768006732399
$ws.Range("K8").Value = 1.030014345216514
$ws.Range("L8").Value = 1.050232017409247
$ws.Range("M8").Value = 1.054171957537409
$ws.Range("N8").Value = 1.029139490896141

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.02036092330407
$ws.Range("D9").Value = 1.025626199846617
$ws.Range("E9").Value = 1.044546865271486
$ws.Range("F9").Value = 1.048318095975707
$ws.Range("I9").Value = 1.028079969387427
$ws.Range("J9").Value = 1.026259913415641
$ws.Range("K9").Value = 1.028822993849234
$ws.Range("L9").Value = 1.047681811851747
$ws.Range("M9").Value = 1.051441007289564
$ws.Range("N9").Value = 1.027717320206349

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.01885739159081
$ws.Range("D10").Value = 1.02453662630092
$ws.Range("E10").Value = 1.042559627613019
$ws.Range("F10").Value = 1.046211500115966
$ws.Range("I10").Value = 1.027873437265552
$ws.Range("J10").Value = 1.025309827006338
$ws.Range("K10").Value = 1.02802368843607
$ws.Range("L10").Value = 1.04598195196852
$ws.Range("M10").Value = 1.049621006598129
$ws.Range("N10").Value = 1.026765884565368

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.018206253128339
$ws.Range("D11").Value = 1.02406446078173
$ws.Range("E11").Value = 1.041699833800765
$ws.Range("F11").Value = 1.045300086051041
$ws.Range("I11").Value = 1.027782175051219
$ws.Range("J11").Value = 1.024897644028609
$ws.Range("K11").Value = 1.027676382752925
$ws.Range("L11").Value = 1.045245934945507
$ws.Range("M11").Value = 1.048833051614521
$ws.Range("N11").Value = 1.026353116240532

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.017964376120332
$ws.Range("D12").Value = 1.02388902211892
$ws.Range("E12").Value = 1.041380569717146
$ws.Range("F12").Value = 1.044961657218132
$ws.Range("I12").Value = 1.027748001153802
$ws.Range("J12").Value = 1.024744422384244
$ws.Range("K12").Value = 1.027547197839387
$ws.Range("L12").Value = 1.044972548289365
$ws.Range("M12").Value = 1.048540385468161
$ws.Range("N12").Value = 1.026199677003858

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.01801626025142
$ws.Range("D13").Value = 1.023926656802348
$ws.Range("E13").Value = 1.041449048388869
$ws.Range("F13").Value = 1.045034246373892
$ws.Range("I13").Value = 1.027755344013186
$ws.Range("J13").Value = 1.024777294310253
$ws.Range("K13").Value = 1.027574916587181
$ws.Range("L13").Value = 1.045031190562054
$ws.Range("M13").Value = 1.04860316268779
$ws.Range("N13").Value = 1.026232595611772

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.01818625981087
$ws.Range("D14").Value = 1.024049960092723
$ws.Range("E14").Value = 1.041673441273597
$ws.Range("F14").Value = 1.045272109170354
$ws.Range("I14").Value = 1.027779355835607
$ws.Range("J14").Value = 1.024884981097476
$ws.Range("K14").Value = 1.027665707960081
$ws.Range("L14").Value = 1.045223336674062
$ws.Range("M14").Value = 1.04880885943781
$ws.Range("N14").Value = 1.026340435326585

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.018291000129042
$ws.Range("D15").Value = 1.024125924002898
$ws.Range("E15").Value = 1.041811710525936
$ws.Range("F15").Value = 1.04541867891069
$ws.Range("I15").Value = 1.027794113873062
$ws.Range("J15").Value = 1.02495131477569
$ws.Range("K15").Value = 1.02772162366407
$ws.Range("L15").Value = 1.045341724545501
$ws.Range("M15").Value = 1.048935597971357
$ws.Range("N15").Value = 1.026406863206231

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.018900603390662
$ws.Range("D16").Value = 1.024567954570446
$ws.Range("E16").Value = 1.042616703686149
$ws.Range("F16").Value = 1.046272003366446
$ws.Range("I16").Value = 1.027879455452593
$ws.Range("J16").Value = 1.025337165607779
$ws.Range("K16").Value = 1.028046712692596
$ws.Range("L16").Value = 1.046030799512903
$ws.Range("M16").Value = 1.049673302858718
$ws.Range("N16").Value = 1.026793261990759

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.019282964632765
$ws.Range("D17").Value = 1.02484512928444
$ws.Range("E17").Value = 1.043121838436676
$ws.Range("F17").Value = 1.046807471761122
$ws.Range("I17").Value = 1.027932497555173
$ws.Range("J17").Value = 1.025578988359862
$ws.Range("K17").Value = 1.028250311108713
$ws.Range("L17").Value = 1.046463045458004
$ws.Range("M17").Value = 1.050136075497481
$ws.Range("N17").Value = 1.027035428158888

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.019505979912666
$ws.Range("D18").Value = 1.025006764573513
$ws.Range("E18").Value = 1.043416542044375
$ws.Range("F18").Value = 1.047119874520397
$ws.Range("I18").Value = 1.02796325929495
$ws.Range("J18").Value = 1.025719963321444
$ws.Range("K18").Value = 1.028368950702512
$ws.Range("L18").Value = 1.046715170586795
$ws.Range("M18").Value = 1.050406014416553
$ws.Range("N18").Value = 1.027176603321085

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.01958202074595
$ws.Range("D19").Value = 1.025061871897875
$ws.Range("E19").Value = 1.043517039809122
$ws.Range("F19").Value = 1.047226408267678
$ws.Range("I19").Value = 1.027973718253492
$ws.Range("J19").Value = 1.025768019220923
$ws.Range("K19").Value = 1.028409384065945
$ws.Range("L19").Value = 1.046801139380773
$ws.Range("M19").Value = 1.050498058624534
$ws.Range("N19").Value = 1.027224727465453

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.019241941881136
$ws.Range("D20").Value = 1.024815394782827
$ws.Range("E20").Value = 1.043067635340741
$ws.Range("F20").Value = 1.046750013527152
$ws.Range("I20").Value = 1.027926824927119
$ws.Range("J20").Value = 1.02555305094412
$ws.Range("K20").Value = 1.028228478912854
$ws.Range("L20").Value = 1.046416669224632
$ws.Range("M20").Value = 1.05008642320243
$ws.Range("N20").Value = 1.027009453909039

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.018136199599272
$ws.Range("D21").Value = 1.024013651869971
$ws.Range("E21").Value = 1.041607360376953
$ws.Range("F21").Value = 1.045202061464204
$ws.Range("I21").Value = 1.02777229254299
$ws.Range("J21").Value = 1.024853273294406
$ws.Range("K21").Value = 1.027638977133114
$ws.Range("L21").Value = 1.045166754372526
$ws.Range("M21").Value = 1.048748286463904
$ws.Range("N21").Value = 1.026308682494798

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.017440887513628
$ws.Range("D22").Value = 1.023509244277686
$ws.Range("E22").Value = 1.040689814243765
$ws.Range("F22").Value = 1.044229443079817
$ws.Range("I22").Value = 1.027673540685237
$ws.Range("J22").Value = 1.024412609532611
$ws.Range("K22").Value = 1.027267292020311
$ws.Range("L22").Value = 1.044380898623717
$ws.Range("M22").Value = 1.047907034897834
$ws.Range("N22").Value = 1.025867392939916

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.017809494040823
$ws.Range("D23").Value = 1.023776670396219
$ws.Range("E23").Value = 1.041176168041377
$ws.Range("F23").Value = 1.044744986651435
$ws.Range("I23").Value = 1.027726041657036
$ws.Range("J23").Value = 1.02464627871904
$ws.Range("K23").Value = 1.027464427982425
$ws.Range("L23").Value = 1.044797494906605
$ws.Range("M23").Value = 1.048352990655702
$ws.Range("N23").Value = 1.026101393963392

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.019260478326146
$ws.Range("D24").Value = 1.024828830635521
$ws.Range("E24").Value = 1.043092127178825
$ws.Range("F24").Value = 1.046775976203272
$ws.Range("I24").Value = 1.027929388689898
$ws.Range("J24").Value = 1.025564771181037
$ws.Range("K24").Value = 1.028238344300518
$ws.Range("L24").Value = 1.046437624637232
$ws.Range("M24").Value = 1.050108858901287
$ws.Range("N24").Value = 1.027021190790038

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.020943890509928
$ws.Range("D25").Value = 1.026048388017325
$ws.Range("E25").Value = 1.045318130411283
$ws.Range("F25").Value = 1.049135703511504
$ws.Range("I25").Value = 1.028158424283765
$ws.Range("J25").Value = 1.026627642734657
$ws.Range("K25").Value = 1.029131882550441
$ws.Range("L25").Value = 1.048341044320153
$ws.Range("M25").Value = 1.05214690392385
$ws.Range("N25").Value = 1.028085571743176
